$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.885.54"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +5.50%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.232.18"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.66%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.74"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.628"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.72"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.08%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.82%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0899"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +5.70%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.561.74"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.67"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.07"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.804"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.60"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.251.04"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +3.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.767.29"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +5.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0903"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.10"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.05"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.69%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "250.35"
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.38"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.07%  "
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.40"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.68"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.29%  "
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.142"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.46%  "
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "167.65"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.86%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.03"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.96%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.69"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.03"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +7.19%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0639"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.45%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -4.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.66"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -5.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.38"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.93%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +30.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.88"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.61%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +4.39%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.60"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +8.88%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.23"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.85%  "
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0978"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +6.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "99.20"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.481.32"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.61"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -6.13%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "52.50"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +6.29%  "
